$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new login-data row (A3: email, B3: password/pin)
$ws.Range("A3").Value = "kalpeshk354@gmail.com"
$ws.Range("B3").Value = 123456

# Turn A3 into a mailto hyperlink pointing at the same address
[void]$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:kalpeshk354@gmail.com")

# Hyperlinks.Add auto-applies the built-in "Hyperlink" cell style (blue/underline);
# reset the cell back to Normal so the data cell keeps its original formatting.
$ws.Range("A3").Style = "Normal"

# Move the active selection to B4 (below the newly added row)
[void]$ws.Range("B4").Select()
